$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Locality"
$ws.Range("B2").Value = 3387.042544212517
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "Significant"

$ws.Range("A3").Value = "Type of property"
$ws.Range("B3").Value = 230.8120164564934
$ws.Range("C3").Value = [double]"3.96512301509665e-52"
$ws.Range("D3").Value = "Significant"

$ws.Range("A4").Value = "Subtype of property"
$ws.Range("B4").Value = 1275.10090185504
$ws.Range("C4").Value = [double]"5.133558478747024e-257"
$ws.Range("D4").Value = "Significant"

$ws.Range("A5").Value = "State of the building"
$ws.Range("B5").Value = 572.6279869382911
$ws.Range("C5").Value = [double]"1.86681536684807e-120"
$ws.Range("D5").Value = "Significant"
